$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40 (ALC)
$ws.Range("H40").Value = 1651.2858
$ws.Range("I40").Value = 900
$ws.Range("J40").Value = 1951.8
$ws.Range("K40").Value = 900
$ws.Range("L40").Value = 1951.8
$ws.Range("M40").Value = -725
$ws.Range("N40").Value = -2301.8

# Row 80 (ALC)
$ws.Range("H80").Value = 859.25
$ws.Range("I80").Value = 859.25
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2577.75
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1579.75
$ws.Range("N80").ClearContents()

# Row 81 (ALC)
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# Row 83 (ALC)
$ws.Range("H83").Value = 859.25
$ws.Range("I83").Value = 859.25
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 7733.25
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -2741.25
$ws.Range("N83").ClearContents()

# Row 84 (ALC)
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# Row 98 (ALC)
$ws.Range("H98").Value = 679.1818
$ws.Range("I98").Value = 679.1818
$ws.Range("K98").Value = 679.1818
$ws.Range("M98").Value = 818.8182

# Row 101 (ALC)
$ws.Range("H101").Value = 800.8333
$ws.Range("I101").Value = 800.8333
$ws.Range("K101").Value = 2402.4999
$ws.Range("M101").Value = -780.4998999999998

# Row 103 (ALC)
$ws.Range("H103").Value = 936.25
$ws.Range("J103").Value = 965
$ws.Range("L103").Value = 2895
$ws.Range("N103").Value = -4067

# Row 122 (ALC)
$ws.Range("H122").Value = 679.1818
$ws.Range("I122").Value = 679.1818
$ws.Range("K122").Value = 2037.5454
$ws.Range("M122").Value = 412.4546

# Row 135 (ALC)
$ws.Range("H135").Value = 1019.9091
$ws.Range("I135").Value = 1022.7778
$ws.Range("K135").Value = 9205.0002
$ws.Range("M135").Value = -6670.0002

# Row 137 (ALC)
$ws.Range("H137").Value = 3818.9614
$ws.Range("I137").Value = 2153.0667
$ws.Range("K137").Value = 6459.2001
$ws.Range("M137").Value = -3909.2001

# Row 138 (ALC)
$ws.Range("H138").Value = 2714.0356
$ws.Range("I138").Value = 1127.0741
$ws.Range("K138").Value = 3381.2223
$ws.Range("M138").Value = 1758.7777

$ws = $wb.Worksheets.Item("ARM")
# Row 5 (ARM)
$ws.Range("H5").Value = 90.84614999999999
$ws.Range("I5").Value = 93.666664
$ws.Range("J5").Value = 90
$ws.Range("K5").Value = 93.666664
$ws.Range("L5").Value = 90
$ws.Range("M5").Value = 18.333336
$ws.Range("N5").Value = -314

# Row 16 (ARM)
$ws.Range("H16").Value = 15007
$ws.Range("J16").Value = 15007
$ws.Range("L16").Value = 15007
$ws.Range("N16").Value = -15581

# Row 61 (ARM)
$ws.Range("H61").Value = 2535.3125
$ws.Range("I61").Value = 2504.5334
$ws.Range("K61").Value = 2504.5334
$ws.Range("M61").Value = -2292.5334

# Row 101 (ARM)
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490

# Row 102 (ARM)
$ws.Range("H102").Value = 2338.2307
$ws.Range("I102").Value = 1908.0834
$ws.Range("K102").Value = 1908.0834
$ws.Range("M102").Value = -286.0834

# Row 136 (ARM)
$ws.Range("H136").Value = 2535.3125
$ws.Range("I136").Value = 2504.5334
$ws.Range("K136").Value = 7513.600199999999
$ws.Range("M136").Value = -4963.600199999999

$ws = $wb.Worksheets.Item("BSM")
# Row 4 (BSM)
$ws.Range("H4").Value = 90.84614999999999
$ws.Range("I4").Value = 93.666664
$ws.Range("J4").Value = 90
$ws.Range("K4").Value = 93.666664
$ws.Range("L4").Value = 90
$ws.Range("M4").Value = 21.333336
$ws.Range("N4").Value = -320

# Row 16 (BSM)
$ws.Range("H16").Value = 500
$ws.Range("I16").Value = 500
$ws.Range("K16").Value = 500
$ws.Range("M16").Value = -330

# Row 87 (BSM)
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

# Row 90 (BSM)
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

# Row 96 (BSM)
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()

# Row 105 (BSM)
$ws.Range("H105").Value = 2455.7856
$ws.Range("I105").Value = 2333.7273
$ws.Range("J105").Value = 2903.3333
$ws.Range("K105").Value = 2333.7273
$ws.Range("L105").Value = 2903.3333
$ws.Range("M105").Value = -586.7273
$ws.Range("N105").Value = -6397.3333

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 2545.3333
$ws.Range("I31").Value = 2705.7144
$ws.Range("K31").Value = 2705.7144
$ws.Range("M31").Value = -2410.7144

# Row 34 (CRP)
$ws.Range("H34").Value = 2545.3333
$ws.Range("I34").Value = 2705.7144
$ws.Range("K34").Value = 2705.7144
$ws.Range("M34").Value = -2503.7144

# Row 58 (CRP)
$ws.Range("H58").Value = 2555.2144
$ws.Range("I58").Value = 2434.2727
$ws.Range("K58").Value = 2434.2727
$ws.Range("M58").Value = -2231.2727

# Row 94 (CRP)
$ws.Range("I94").Value = 3000
$ws.Range("K94").Value = 3000
$ws.Range("M94").Value = -2549

# Row 132 (CRP)
$ws.Range("H132").Value = 3174.2083
$ws.Range("I132").Value = 2913.524
$ws.Range("K132").Value = 8740.572
$ws.Range("M132").Value = -6210.572

# Row 134 (CRP)
$ws.Range("H134").Value = 4024.6667
$ws.Range("I134").Value = 4027.8125
$ws.Range("J134").Value = 3999.5
$ws.Range("K134").Value = 12083.4375
$ws.Range("L134").Value = 11998.5
$ws.Range("M134").Value = -9548.4375
$ws.Range("N134").Value = -17068.5

# Row 136 (CRP)
$ws.Range("H136").Value = 2555.2144
$ws.Range("I136").Value = 2434.2727
$ws.Range("K136").Value = 7302.8181
$ws.Range("M136").Value = -4752.8181

$ws = $wb.Worksheets.Item("CUL")
# Row 33 (CUL)
$ws.Range("H33").Value = 826.8570999999999
$ws.Range("I33").Value = 131
$ws.Range("K33").Value = 786
$ws.Range("M33").Value = -503

# Row 107 (CUL)
$ws.Range("H107").Value = 1832.6666
$ws.Range("J107").Value = 199.2
$ws.Range("L107").Value = 597.5999999999999
$ws.Range("N107").Value = -4437.6

$ws = $wb.Worksheets.Item("GSM")
# Row 27 (GSM)
$ws.Range("H27").Value = 10011
$ws.Range("J27").Value = 10011
$ws.Range("L27").Value = 10011
$ws.Range("N27").Value = -10343

# Row 126 (GSM)
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

# Row 132 (GSM)
$ws.Range("H132").Value = 2997.6
$ws.Range("I132").Value = 2998.6667
$ws.Range("K132").Value = 8996.000100000001
$ws.Range("M132").Value = -6466.000100000001

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (LTW)
$ws.Range("H22").Value = 2000
$ws.Range("I22").Value = 2000
$ws.Range("K22").Value = 2000
$ws.Range("M22").Value = -1705

# Row 27 (LTW)
$ws.Range("H27").Value = 2000
$ws.Range("I27").Value = 2000
$ws.Range("K27").Value = 2000
$ws.Range("M27").Value = -1893

# Row 43 (LTW)
$ws.Range("H43").Value = 813076.6
$ws.Range("J43").Value = 813076.6
$ws.Range("L43").Value = 813076.6
$ws.Range("N43").Value = -813462.6

# Row 46 (LTW)
$ws.Range("H46").Value = 499.5
$ws.Range("I46").Value = 499.5
$ws.Range("K46").Value = 499.5
$ws.Range("M46").Value = -311.5

# Row 82 (LTW)
$ws.Range("H82").Value = 1434.8
$ws.Range("I82").Value = 1321.2858
$ws.Range("J82").Value = 1699.6666
$ws.Range("K82").Value = 1321.2858
$ws.Range("L82").Value = 1699.6666
$ws.Range("M82").Value = -960.2858000000001
$ws.Range("N82").Value = -2421.6666

# Row 85 (LTW)
$ws.Range("H85").Value = 1434.8
$ws.Range("I85").Value = 1321.2858
$ws.Range("J85").Value = 1699.6666
$ws.Range("K85").Value = 1321.2858
$ws.Range("L85").Value = 1699.6666
$ws.Range("M85").Value = -73.28580000000011
$ws.Range("N85").Value = -4195.6666

# Row 122 (LTW)
$ws.Range("H122").Value = 3398.5
$ws.Range("I122").Value = 2600.25
$ws.Range("K122").Value = 7800.75
$ws.Range("M122").Value = -5350.75

# Row 132 (LTW)
$ws.Range("H132").Value = 1570.7142
$ws.Range("I132").Value = 1249.6666
$ws.Range("J132").Value = 3497
$ws.Range("K132").Value = 3748.9998
$ws.Range("L132").Value = 10491
$ws.Range("M132").Value = -1218.9998
$ws.Range("N132").Value = -15551

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (WVR)
$ws.Range("H81").Value = 4709.4546
$ws.Range("I81").Value = 1088.8889
$ws.Range("K81").Value = 2177.7778
$ws.Range("M81").Value = -1116.7778

# Row 84 (WVR)
$ws.Range("H84").Value = 4709.4546
$ws.Range("I84").Value = 1088.8889
$ws.Range("K84").Value = 10888.889
$ws.Range("M84").Value = -5584.888999999999

# Row 101 (WVR)
$ws.Range("H101").Value = 55000
$ws.Range("J101").Value = 55000
$ws.Range("L101").Value = 55000
$ws.Range("N101").Value = -61490

# Row 132 (WVR)
$ws.Range("H132").Value = 1345.7142
$ws.Range("I132").Value = 1263.1666
$ws.Range("J132").Value = 1841
$ws.Range("K132").Value = 3789.4998
$ws.Range("L132").Value = 5523
$ws.Range("M132").Value = -1259.4998
$ws.Range("N132").Value = -10583

# Row 141 (WVR)
$ws.Range("H141").Value = 120000
$ws.Range("J141").Value = 120000
$ws.Range("L141").Value = 120000
$ws.Range("N141").Value = -130360
